# TC05_INS_Filter_Doc-DCTD.xlsx - "Fixed Bento 80 Test scripts"
#
# The Neo4j query stored in startup!B2 (the "dbExcel" query used to pull the
# DCTD project list) was missing a deterministic ordering/row cap, which made
# the automated test flaky/slow. Append an ORDER BY + LIMIT clause to the end
# of the existing Cypher query text, leaving everything else untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$query = @'
MATCH (p:project)-->(pgm:program)
where p.lead_doc='DCTD'
WITH DISTINCT p, pgm
RETURN
coalesce(p.project_id, '') AS `Project ID`,
coalesce (pgm.program_id, '')AS `Program`,
coalesce(p.project_title, '') AS `Project Title`,
coalesce(p.principal_investigators, '') AS `Principal Investigators`,
coalesce(p.program_officers, '') AS `Program Officers`,
coalesce(p.lead_doc, '')AS `Lead DOC`,
SUBSTRING(p.project_id, 1, 3) AS `Activity code`,
"$" + apoc.number.format(toInteger(p.award_amount)) AS `Award Amount`,
coalesce(p.project_end_date, '') AS `Project End Date`,
coalesce(p.fiscal_year,'')AS `Fiscal Year` ORDER BY p.project_id ASC LIMIT 100
'@

$ws.Range("B2").Value = $query
